$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (previously for MuSCs target cluster; no longer needed in updated TPM data)
$ws.Rows.Item(8).Delete() | Out-Null
$ws.Rows.Item(8).Delete() | Out-Null
$ws.Rows.Item(8).Delete() | Out-Null

# Update remaining rows 2-7 with refreshed TPM-derived values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.630712666666667
$ws.Range("H2").Value = 7.892138000000001
$ws.Range("I2").Value = 0.3947434022685045
$ws.Range("J2").Value = 0.3947434022685045
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.02114766666666667
$ws.Range("N2").Value = 0.063443
$ws.Range("O2").Value = 0.7642906190894964
$ws.Range("P2").Value = 0.7642906190894964
$ws.Range("Q2").Value = 0.05563343457044445
$ws.Range("R2").Value = 0.5007009111340001
$ws.Range("S2").Value = 0.3016986793012894
$ws.Range("T2").Value = 0.3016986793012894
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk3"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.630712666666667
$ws.Range("H3").Value = 7.892138000000001
$ws.Range("I3").Value = 0.3947434022685045
$ws.Range("J3").Value = 0.3947434022685045
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.006522
$ws.Range("N3").Value = 0.019566
$ws.Range("O3").Value = 0.2357093809105037
$ws.Range("P3").Value = 0.2357093809105037
$ws.Range("Q3").Value = 0.017157508012
$ws.Range("R3").Value = 0.154417572108
$ws.Range("S3").Value = 0.09304472296721512
$ws.Range("T3").Value = 0.09304472296721512
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.99311
$ws.Range("H4").Value = 11.97933
$ws.Range("I4").Value = 0.5991736942634763
$ws.Range("J4").Value = 0.5991736942634763
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.02114766666666667
$ws.Range("N4").Value = 0.063443
$ws.Range("O4").Value = 0.7642906190894964
$ws.Range("P4").Value = 0.7642906190894964
$ws.Range("Q4").Value = 0.08444495924333333
$ws.Range("R4").Value = 0.76000463319
$ws.Range("S4").Value = 0.4579428337307729
$ws.Range("T4").Value = 0.4579428337307729
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk3"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.99311
$ws.Range("H5").Value = 11.97933
$ws.Range("I5").Value = 0.5991736942634763
$ws.Range("J5").Value = 0.5991736942634763
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.006522
$ws.Range("N5").Value = 0.019566
$ws.Range("O5").Value = 0.2357093809105037
$ws.Range("P5").Value = 0.2357093809105037
$ws.Range("Q5").Value = 0.02604306342
$ws.Range("R5").Value = 0.23438757078
$ws.Range("S5").Value = 0.1412308605327034
$ws.Range("T5").Value = 0.1412308605327034
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04053866666666667
$ws.Range("H6").Value = 0.121616
$ws.Range("I6").Value = 0.00608290346801924
$ws.Range("J6").Value = 0.006082903468019241
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.02114766666666667
$ws.Range("N6").Value = 0.063443
$ws.Range("O6").Value = 0.7642906190894964
$ws.Range("P6").Value = 0.7642906190894964
$ws.Range("Q6").Value = 0.0008572982097777778
$ws.Range("R6").Value = 0.007715683888
$ws.Range("S6").Value = 0.00464910605743407
$ws.Range("T6").Value = 0.004649106057434071
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04053866666666667
$ws.Range("H7").Value = 0.121616
$ws.Range("I7").Value = 0.00608290346801924
$ws.Range("J7").Value = 0.006082903468019241
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.006522
$ws.Range("N7").Value = 0.019566
$ws.Range("O7").Value = 0.2357093809105037
$ws.Range("P7").Value = 0.2357093809105037
$ws.Range("Q7").Value = 0.000264393184
$ws.Range("R7").Value = 0.002379538656
$ws.Range("S7").Value = 0.001433797410585171
$ws.Range("T7").Value = 0.001433797410585171
